$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A31").Value = "What's the maximum number of lithology types in an log?"
$ws.Range("B31").Value = "llama3.2:latest"
$ws.Range("C31").Value = "The maximum number of lithology types that can be recorded in a log is 450."
